# Rename sheets and update the "intensity" (row 2) and "energy" (row 5)
# values on each of the 8 RSA benchmark worksheets to match the new,
# lower-precision "compare tab" numbers.

$wb = $excel.ActiveWorkbook

$sheetData = @(
    @{ Name = "data_CCM code_FLASH_enc"; B2 = 12457; C2 = 21646; D2 = 29400; B5 = 0.672; C5 = 0.656; D5 = 0.7   },
    @{ Name = "data_CCM code_FLASH_dec"; B2 = 12442; C2 = 22461; D2 = 30252; B5 = 5.212; C5 = 5.251; D5 = 5.549 },
    @{ Name = "data_CCM code_CCM_enc";   B2 = 11997; C2 = 21316; D2 = 29356; B5 = 0.672; C5 = 0.665; D5 = 0.696 },
    @{ Name = "data_CCM code_CCM_dec";   B2 = 12020; C2 = 22149; D2 = 30617; B5 = 5.245; C5 = 5.261; D5 = 5.428 },
    @{ Name = "data_RAM code_FLASH_enc"; B2 = 12464; C2 = 22098; D2 = 29422; B5 = 0.666; C5 = 0.665; D5 = 0.701 },
    @{ Name = "data_RAM code_FLASH_dec"; B2 = 12501; C2 = 22094; D2 = 30260; B5 = 5.229; C5 = 5.171; D5 = 5.56  },
    @{ Name = "data_RAM code_CCM_enc";   B2 = 11980; C2 = 21350; D2 = 29292; B5 = 0.668; C5 = 0.664; D5 = 0.694 },
    @{ Name = "data_RAM code_CCM_dec";   B2 = 11997; C2 = 21781; D2 = 30817; B5 = 5.239; C5 = 5.167; D5 = 5.469 }
)

for ($i = 0; $i -lt $sheetData.Count; $i++) {
    $info = $sheetData[$i]
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $info.Name

    $ws.Range("B2").Value = $info.B2
    $ws.Range("C2").Value = $info.C2
    $ws.Range("D2").Value = $info.D2

    $ws.Range("B5").Value = $info.B5
    $ws.Range("C5").Value = $info.C5
    $ws.Range("D5").Value = $info.D5
}
